# Generate Report for Handoff
# b.md has been handed off for localization; update its status, the newly
# generated handoff xliff file names/timestamps, and the version-mismatch
# error detail on the zh-cn and de-de status sheets, plus the Overview
# summary sheet.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/28d2d1ca24872018e84eaf69837a977deb460028/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/52a55455522cbfd75fef5458eba567336a58621c/e2e/b.md."

# ---- Overview sheet: b.md row (row 3) now "Ready for handoff" ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-05 02:42:39"

# ---- zh-cn sheet: b.md row (row 3) ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-05 02:42:34"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.1667

# ---- de-de sheet: b.md row (row 3) ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
# "False" would otherwise auto-coerce to a native boolean cell; force it to
# stay text (matching the source workbook's shared-string "False") via a
# leading apostrophe, then restore the default (non quote-prefixed) style.
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("F3").Style = "Normal"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-05 02:42:39"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.1667
